$d = $word.ActiveDocument

# The document starts with two consecutive "Title" style paragraphs: the
# first contains the stray text "Test" (followed by the _GoBack
# bookmark), the second is completely empty. The edit deletes the word
# "Test" together with the paragraph mark right after it, which merges
# the two paragraphs into a single empty (bookmark-only) Title
# paragraph - exactly what's left once the literal text and the extra
# paragraph break both disappear.
#
# "^p" (outside of wildcard mode) is Word's special-character code for
# a paragraph mark, so searching for "Test^p" matches the word plus the
# paragraph break that follows it; replacing with "" removes both in a
# single operation while leaving the bookmark and the paragraph's
# formatting (Title style, line spacing) intact.
$found = $d.Content.Find.Execute(
    "Test^p", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

$d.Save()
